$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "Big-O Guide" slide at position 12 (Title and Content
#    layout), pushing "Space Complexity" and everything after it down by one.
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add(12, 2)

$newTitle = $newSlide.Shapes.Item(1)
$newTitle.TextFrame.TextRange.Text = "Big-O Guide "

$newBody = $newSlide.Shapes.Item(2)
$bodyRange = $newBody.TextFrame.TextRange
$bodyRange.Text = "Calculations not dependent on the input size – O(1)"
$bodyRange.InsertAfter("`rLoop- O(n)") | Out-Null
$bodyRange.InsertAfter("`rNested loops – O(n^2)") | Out-Null
$bodyRange.InsertAfter("`rInput size reduced by half – O(logn)") | Out-Null
$bodyRange.InsertAfter("`r") | Out-Null

# ---------------------------------------------------------------------------
# 2. "Fibonacci sequence" slide (now at position 20) gets one more example
#    appended to its bullet list.
# ---------------------------------------------------------------------------
$fibSlide = $p.Slides.Item(20)
$fibBody = $fibSlide.Shapes.Item(2)
$fibRange = $fibBody.TextFrame.TextRange
$fibRange.InsertAfter("`rFibonacci of (7) = [0,1,1,2,3,5,8] every number is sum of previous two numbers") | Out-Null

# ---------------------------------------------------------------------------
# 3. The blank slide at position 21 becomes "Factorial of a number".
# ---------------------------------------------------------------------------
$factSlide = $p.Slides.Item(21)
$factTitle = $factSlide.Shapes.Item(1)
$factTitle.TextFrame.TextRange.Text = "Factorial of a number"

$factBody = $factSlide.Shapes.Item(2)
$factRange = $factBody.TextFrame.TextRange
$factRange.Text = "Problem – give an integer ‘n’, find the factorial of that integer."
$factRange.InsertAfter("`rIn mathematics, the factorial of a non-negative integer ‘n’, denoted n!, is the product of all positive integers less than or equal to ‘n’.\") | Out-Null
$factRange.InsertAfter("`rFactorial of zero is 1.") | Out-Null
$factRange.InsertAfter("`rFactorial of (4) = 4*3*2*1 =24") | Out-Null
$factRange.InsertAfter("`rFactorial(5) = 5*4*3*2*1 = 120") | Out-Null

# ---------------------------------------------------------------------------
# 4. The blank slide at position 22 becomes "Prime Number".
# ---------------------------------------------------------------------------
$primeSlide = $p.Slides.Item(22)
$primeTitle = $primeSlide.Shapes.Item(1)
$primeTitle.TextFrame.TextRange.Text = "Prime Number"

$primeBody = $primeSlide.Shapes.Item(2)
$primeRange = $primeBody.TextFrame.TextRange
$primeRange.Text = "Problem – Give a natural number ‘n’, determine  if the number is prime or not."
